# Add the new "2022-Q3" quarterly sheet.
#
# The workbook currently has:
#   总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3
#
# After this script it should have:
#   总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3
#
# "2022-Q3" is inserted right after "总计" (i.e. in front of the other
# quarter sheets), built as a copy of the "2022-Q2" sheet (so it keeps the
# same layout/formatting/column headers) with refreshed figures. The
# "总计" (summary) sheet gets a new row at the top for 2022-Q3 and all
# the other rows shift down by one.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "2022-Q2" sheet, placing the copy right before it,
#        then rename the copy to "2022-Q3". This preserves all existing
#        formatting/styles/column widths from the template sheet. ---
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# --- 2. Update the figures on the new "2022-Q3" sheet. ---
# Columns D-G on these sheets are stored as text, not numbers, so force a
# text number-format before writing so the values keep their original type.
$dataRange = $q3.Range("D2:G3")
$dataRange.NumberFormat = "@"

$q3.Range("D2").Value = "0.26"
$q3.Range("E2").Value = "86.78"
$q3.Range("F2").Value = "6.13"
$q3.Range("G2").Value = "0.0159"

$q3.Range("D3").Value = "0.20"
$q3.Range("E3").Value = "86.78"
$q3.Range("F3").Value = "6.13"
$q3.Range("G3").Value = "0.0123"

# --- 3. Update the "总计" summary sheet: insert a new row for 2022-Q3
#        and shift the existing quarters down by one row. ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.03

# Re-sequence the running index in column A for the rows that moved down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
